$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-07 Wednesday", "2025-05-08 Thursday"),
    @("192×3=576", "501×4=2004"),
    @("902×8=7216", "255×3=765"),
    @("320×8=2560", "827×8=6616"),
    @("433×8=3464", "698×3=2094"),
    @("842×5=4210", "829×8=6632"),
    @("772×2=1544", "286×8=2288"),
    @("772×3=2316", "410×9=3690"),
    @("801×9=7209", "213×8=1704"),
    @("933×5=4665", "105×8=840"),
    @("640×9=5760", "634×3=1902"),
    @("514×7=3598", "389×5=1945"),
    @("401×8=3208", "394×3=1182"),
    @("575×2=1150", "425×9=3825"),
    @("206×5=1030", "632×9=5688"),
    @("149×8=1192", "758×4=3032"),
    @("772×4=3088", "778×3=2334"),
    @("396×6=2376", "339×9=3051"),
    @("855×8=6840", "902×9=8118"),
    @("151×9=1359", "535×3=1605"),
    @("389×3=1167", "314×8=2512"),
    @("156×5=780", "928×8=7424"),
    @("196×7=1372", "956×8=7648"),
    @("455×7=3185", "721×2=1442"),
    @("203×6=1218", "444×3=1332"),
    @("123×2=246", "417×6=2502")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
